# Auto-generated edit script: updates cached price/profit figures on the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# freshly-pulled marketboard prices from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 138.9
$ws.Cells.Item(11, 9).Value = 138.9
$ws.Cells.Item(11, 11).Value = 138.9
$ws.Cells.Item(11, 13).Value = 1.099999999999994
$ws.Cells.Item(19, 8).Value = 1956
$ws.Cells.Item(19, 9).Value = 1720.2
$ws.Cells.Item(19, 11).Value = 1720.2
$ws.Cells.Item(19, 13).Value = -1545.2
$ws.Cells.Item(98, 8).Value = 843.5
$ws.Cells.Item(98, 9).Value = 843.5
$ws.Cells.Item(98, 11).Value = 843.5
$ws.Cells.Item(98, 13).Value = 654.5
$ws.Cells.Item(111, 8).Value = 1800.5
$ws.Cells.Item(111, 9).Value = 1660
$ws.Cells.Item(111, 10).Value = 2222
$ws.Cells.Item(111, 11).Value = 4980
$ws.Cells.Item(111, 12).Value = 6666
$ws.Cells.Item(111, 13).Value = -1913
$ws.Cells.Item(111, 14).Value = -12800
$ws.Cells.Item(112, 8).Value = 1106.04
$ws.Cells.Item(112, 10).Value = 1391.2354
$ws.Cells.Item(112, 12).Value = 4173.706200000001
$ws.Cells.Item(112, 14).Value = -6389.706200000001
$ws.Cells.Item(116, 8).Value = 3906.2307
$ws.Cells.Item(116, 9).Value = 4076.5557
$ws.Cells.Item(116, 10).Value = 3523
$ws.Cells.Item(116, 11).Value = 4076.5557
$ws.Cells.Item(116, 12).Value = 3523
$ws.Cells.Item(116, 13).Value = -634.5556999999999
$ws.Cells.Item(116, 14).Value = -10407
$ws.Cells.Item(122, 8).Value = 843.5
$ws.Cells.Item(122, 9).Value = 843.5
$ws.Cells.Item(122, 11).Value = 2530.5
$ws.Cells.Item(122, 13).Value = -80.5
$ws.Cells.Item(132, 8).Value = 2204
$ws.Cells.Item(132, 9).Value = 2204
$ws.Cells.Item(132, 11).Value = 6612
$ws.Cells.Item(132, 13).Value = -4082

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 1600
$ws.Cells.Item(16, 9).Value = 1600
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1600
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -1313
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(45, 8).Value = 2775
$ws.Cells.Item(45, 9).Value = 2775
$ws.Cells.Item(45, 11).Value = 2775
$ws.Cells.Item(45, 13).Value = -2398
$ws.Cells.Item(61, 8).Value = 4260.615
$ws.Cells.Item(61, 9).Value = 4365.9165
$ws.Cells.Item(61, 11).Value = 4365.9165
$ws.Cells.Item(61, 13).Value = -4153.9165
$ws.Cells.Item(110, 8).Value = 2109.7144
$ws.Cells.Item(110, 9).Value = 2294.8333
$ws.Cells.Item(110, 10).Value = 999
$ws.Cells.Item(110, 11).Value = 2294.8333
$ws.Cells.Item(110, 12).Value = 999
$ws.Cells.Item(110, 13).Value = -249.8332999999998
$ws.Cells.Item(110, 14).Value = -5089
$ws.Cells.Item(132, 8).Value = 2099.1538
$ws.Cells.Item(132, 9).Value = 1728.6364
$ws.Cells.Item(132, 11).Value = 5185.9092
$ws.Cells.Item(132, 13).Value = -2655.9092
$ws.Cells.Item(136, 8).Value = 4260.615
$ws.Cells.Item(136, 9).Value = 4365.9165
$ws.Cells.Item(136, 11).Value = 13097.7495
$ws.Cells.Item(136, 13).Value = -10547.7495

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4999.5
$ws.Cells.Item(86, 9).Value = 4999
$ws.Cells.Item(86, 10).Value = 5000
$ws.Cells.Item(86, 11).Value = 4999
$ws.Cells.Item(86, 12).Value = 5000
$ws.Cells.Item(86, 13).Value = -3876
$ws.Cells.Item(86, 14).Value = -7246
$ws.Cells.Item(89, 8).Value = 4999.5
$ws.Cells.Item(89, 9).Value = 4999
$ws.Cells.Item(89, 10).Value = 5000
$ws.Cells.Item(89, 11).Value = 24995
$ws.Cells.Item(89, 12).Value = 25000
$ws.Cells.Item(89, 13).Value = -19379
$ws.Cells.Item(89, 14).Value = -36232

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2989
$ws.Cells.Item(16, 9).Value = 2989.5
$ws.Cells.Item(16, 10).Value = 2988
$ws.Cells.Item(16, 11).Value = 2989.5
$ws.Cells.Item(16, 12).Value = 2988
$ws.Cells.Item(16, 13).Value = -2702.5
$ws.Cells.Item(16, 14).Value = -3562
$ws.Cells.Item(58, 8).Value = 2414.6667
$ws.Cells.Item(58, 10).Value = 2398.5
$ws.Cells.Item(58, 12).Value = 2398.5
$ws.Cells.Item(58, 14).Value = -2804.5
$ws.Cells.Item(74, 8).Value = 34989.668
$ws.Cells.Item(74, 10).Value = 34989.668
$ws.Cells.Item(74, 12).Value = 34989.668
$ws.Cells.Item(74, 14).Value = -36737.668
$ws.Cells.Item(77, 8).Value = 34989.668
$ws.Cells.Item(77, 10).Value = 34989.668
$ws.Cells.Item(77, 12).Value = 104969.004
$ws.Cells.Item(77, 14).Value = -113705.004
$ws.Cells.Item(94, 8).Value = 2999.5
$ws.Cells.Item(94, 10).Value = 2998
$ws.Cells.Item(94, 12).Value = 2998
$ws.Cells.Item(94, 14).Value = -3900
$ws.Cells.Item(113, 8).Value = 2989
$ws.Cells.Item(113, 9).Value = 2989.5
$ws.Cells.Item(113, 10).Value = 2988
$ws.Cells.Item(113, 11).Value = 2989.5
$ws.Cells.Item(113, 12).Value = 2988
$ws.Cells.Item(113, 13).Value = -819.5
$ws.Cells.Item(113, 14).Value = -7328
$ws.Cells.Item(134, 8).Value = 3649.8235
$ws.Cells.Item(134, 9).Value = 3635.077
$ws.Cells.Item(134, 11).Value = 10905.231
$ws.Cells.Item(134, 13).Value = -8370.231
$ws.Cells.Item(136, 8).Value = 2414.6667
$ws.Cells.Item(136, 10).Value = 2398.5
$ws.Cells.Item(136, 12).Value = 7195.5
$ws.Cells.Item(136, 14).Value = -12295.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 59349810
$ws.Cells.Item(4, 9).Value = 5002828
$ws.Cells.Item(4, 11).Value = 15008484
$ws.Cells.Item(4, 13).Value = -15008372
$ws.Cells.Item(12, 8).Value = 345.3
$ws.Cells.Item(12, 9).Value = 281.2857
$ws.Cells.Item(12, 10).Value = 494.66666
$ws.Cells.Item(12, 11).Value = 843.8571000000001
$ws.Cells.Item(12, 12).Value = 1483.99998
$ws.Cells.Item(12, 13).Value = -670.8571000000001
$ws.Cells.Item(12, 14).Value = -1829.99998
$ws.Cells.Item(16, 8).Value = 300
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 13).ClearContents()
$ws.Cells.Item(113, 8).Value = 2796
$ws.Cells.Item(113, 9).Value = 999
$ws.Cells.Item(113, 11).Value = 2997
$ws.Cells.Item(113, 13).Value = -827

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1599.6666
$ws.Cells.Item(97, 9).Value = 2149.5
$ws.Cells.Item(97, 10).Value = 500
$ws.Cells.Item(97, 11).Value = 2149.5
$ws.Cells.Item(97, 12).Value = 500
$ws.Cells.Item(97, 13).Value = -1653.5
$ws.Cells.Item(97, 14).Value = -1492
$ws.Cells.Item(122, 8).Value = 929.1667
$ws.Cells.Item(122, 10).Value = 949.5
$ws.Cells.Item(122, 12).Value = 2848.5
$ws.Cells.Item(122, 14).Value = -7748.5
$ws.Cells.Item(126, 8).Value = 3997.5
$ws.Cells.Item(126, 9).Value = 3997.5
$ws.Cells.Item(126, 11).Value = 11992.5
$ws.Cells.Item(126, 13).Value = -9522.5
$ws.Cells.Item(127, 8).Value = 66666
$ws.Cells.Item(127, 10).Value = 66666
$ws.Cells.Item(127, 12).Value = 66666
$ws.Cells.Item(127, 14).Value = -76586

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(46, 8).Value = 4032.8
$ws.Cells.Item(46, 9).Value = 1999.2858
$ws.Cells.Item(46, 10).Value = 5812.125
$ws.Cells.Item(46, 11).Value = 1999.2858
$ws.Cells.Item(46, 12).Value = 5812.125
$ws.Cells.Item(46, 13).Value = -1811.2858
$ws.Cells.Item(46, 14).Value = -6188.125
$ws.Cells.Item(55, 8).Value = 186.58333
$ws.Cells.Item(55, 9).Value = 140.2
$ws.Cells.Item(55, 10).Value = 219.71428
$ws.Cells.Item(55, 11).Value = 140.2
$ws.Cells.Item(55, 12).Value = 219.71428
$ws.Cells.Item(55, 13).Value = 32.80000000000001
$ws.Cells.Item(55, 14).Value = -565.71428
$ws.Cells.Item(61, 8).Value = 3255.5
$ws.Cells.Item(61, 10).Value = 3150
$ws.Cells.Item(61, 12).Value = 3150
$ws.Cells.Item(61, 14).Value = -3554
$ws.Cells.Item(82, 8).Value = 1630.7778
$ws.Cells.Item(82, 9).Value = 1695.6
$ws.Cells.Item(82, 10).Value = 1549.75
$ws.Cells.Item(82, 11).Value = 1695.6
$ws.Cells.Item(82, 12).Value = 1549.75
$ws.Cells.Item(82, 13).Value = -1334.6
$ws.Cells.Item(82, 14).Value = -2271.75
$ws.Cells.Item(85, 8).Value = 1630.7778
$ws.Cells.Item(85, 9).Value = 1695.6
$ws.Cells.Item(85, 10).Value = 1549.75
$ws.Cells.Item(85, 11).Value = 1695.6
$ws.Cells.Item(85, 12).Value = 1549.75
$ws.Cells.Item(85, 13).Value = -447.5999999999999
$ws.Cells.Item(85, 14).Value = -4045.75
$ws.Cells.Item(100, 8).Value = 1865
$ws.Cells.Item(100, 9).Value = 1838
$ws.Cells.Item(100, 11).Value = 1838
$ws.Cells.Item(100, 13).Value = -1297
$ws.Cells.Item(113, 8).Value = 3255.5
$ws.Cells.Item(113, 10).Value = 3150
$ws.Cells.Item(113, 12).Value = 3150
$ws.Cells.Item(113, 14).Value = -7490

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(15, 8).Value = 14330
$ws.Cells.Item(15, 10).Value = 14330
$ws.Cells.Item(15, 12).Value = 14330
$ws.Cells.Item(15, 14).Value = -14906
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(62, 8).Value = 11299.2
$ws.Cells.Item(62, 9).Value = 8124.25
$ws.Cells.Item(62, 11).Value = 8124.25
$ws.Cells.Item(62, 13).Value = -7500.25
$ws.Cells.Item(65, 8).Value = 11299.2
$ws.Cells.Item(65, 9).Value = 8124.25
$ws.Cells.Item(65, 11).Value = 40621.25
$ws.Cells.Item(65, 13).Value = -37501.25
$ws.Cells.Item(96, 8).Value = 495.33334
$ws.Cells.Item(96, 9).Value = 495.33334
$ws.Cells.Item(96, 11).Value = 495.33334
$ws.Cells.Item(96, 13).Value = 877.66666
$ws.Cells.Item(100, 8).Value = 1126.3334
$ws.Cells.Item(100, 9).Value = 1344.5
$ws.Cells.Item(100, 11).Value = 2689
$ws.Cells.Item(100, 13).Value = -2148
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).ClearContents()
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2379.1155
$ws.Cells.Item(132, 10).Value = 3698
$ws.Cells.Item(132, 12).Value = 11094
$ws.Cells.Item(132, 14).Value = -16154
$ws.Cells.Item(136, 8).Value = 798.1579
$ws.Cells.Item(136, 9).Value = 598.0625
$ws.Cells.Item(136, 11).Value = 1794.1875
$ws.Cells.Item(136, 13).Value = 755.8125
